# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" data table with newer figures and
# re-sorts a handful of rows whose relative ranking (by "Casos totales",
# column B, descending) changed as a result. Because the sheet is always
# kept sorted by column B, a handful of adjacent-row pairs needed to
# swap their country name + figures so the table stays sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 10:05"

# Rusia (row 5) - updated figures
$ws.Range("B5").Value = 317554
$ws.Range("C5").Value = 8849
$ws.Range("D5").Value = 92681
$ws.Range("E5").Value = 221774
$ws.Range("G5").Value = 127
$ws.Range("H5").Value = 3099

# Singapur overtakes Portugal (rows 29-30)
$ws.Range("A29").Value = "Singapur"
$ws.Range("B29").Value = 29812
$ws.Range("C29").Value = 448
$ws.Range("D29").Value = 11207
$ws.Range("E29").Value = 18583
$ws.Range("H29").Value = 22

$ws.Range("A30").Value = "Portugal"
$ws.Range("B30").Value = 29660
$ws.Range("D30").Value = 6452
$ws.Range("E30").Value = 21945
$ws.Range("H30").Value = 1263

# Afganistan overtakes Noruega (rows 53-54)
$ws.Range("A53").Value = "Afganistan"
$ws.Range("B53").Value = 8676
$ws.Range("C53").Value = 531
$ws.Range("D53").Value = 938
$ws.Range("E53").Value = 7545
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 193

$ws.Range("A54").Value = "Noruega"
$ws.Range("B54").Value = 8281
$ws.Range("D54").Value = 32
$ws.Range("E54").Value = 8015
$ws.Range("H54").Value = 234

# Estonia (row 90) - updated figures
$ws.Range("B90").Value = 1800
$ws.Range("C90").Value = 6
$ws.Range("D90").Value = 963
$ws.Range("E90").Value = 773

# Lituania (row 93) - updated figures
$ws.Range("B93").Value = 1594
$ws.Range("C93").Value = 17
$ws.Range("E93").Value = 484
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 61

# Letonia (row 97) - updated figures
$ws.Range("B97").Value = 1502
$ws.Range("C97").Value = 6
$ws.Range("D97").Value = 1245
$ws.Range("E97").Value = 229

# Sri Lanka (row 106) - updated figures
$ws.Range("D106").Value = 604
$ws.Range("E106").Value = 415

# Nueva Caledonia overtakes Santa Lucia (rows 197-198, tie on column B)
$ws.Range("A197").Value = "Nueva Caledonia"
$ws.Range("A198").Value = "Santa Lucia"

# Seychelles overtakes Montserrat (rows 209-210)
$ws.Range("A209").Value = "Seychelles"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

# Sahara Occidental overtakes Bonaire, San Eustaquio y Saba (rows 214-215, tie on column B)
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
